$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new rows (120, 121) following the same pattern as rows 118/119
# (Arabic postal-code location rows), with incrementing codes 10113, 10114.
$newRows = @(
    @(10113, 10113, 5, "الرمز البريدي", "BNMR", "ara", $true, "superadmin", "now()"),
    @(10114, 10114, 5, "الرمز البريدي", "BNMR", "ara", $true, "superadmin", "now()")
)

$startRow = 120
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Update the selection to reflect the post-edit active cell / selection
# (row 122, full-row-to-end selection) as recorded in the saved file.
$ws.Range("A122:XFD1048576").Select()

$wb.Save()
